# ==========================================================================
# Edit script: restructure PlayerPerformance workbook
#   1. Add a new "Player Info" sheet (first position) with player bio data
#   2. On "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, replace the
#      full howstat URL with just the numeric match code, and drop the
#      (empty) INNING_NUMBER placeholder cells that had no value.
#   3. On "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE and replace
#      the full howstat URL with just the numeric match code.
#   4. Add a new "ODI Batting Extra" sheet (last position) with additional
#      per-match batting detail.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# 1. Insert "Player Info" sheet before the current first sheet
# --------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet, $null)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = New-Object 'object[,]' 1,4
$playerInfoHeaders[0,0] = "ID"
$playerInfoHeaders[0,1] = "NAME"
$playerInfoHeaders[0,2] = "BATTING_HAND"
$playerInfoHeaders[0,3] = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Value = $playerInfoHeaders
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

$playerInfo.Range("A2:D2").NumberFormat = "@"
$playerInfoData = New-Object 'object[,]' 1,4
$playerInfoData[0,0] = "3642"
$playerInfoData[0,1] = "Wayne Dillon Parnell"
$playerInfoData[0,2] = "Left Handed"
$playerInfoData[0,3] = "Left Arm Medium Fast"
$playerInfo.Range("A2:D2").Value = $playerInfoData

$playerInfo.Range("A1").Select()

# --------------------------------------------------------------------
# 2. "ODI Batting" sheet updates
# --------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

$batting.Range("D1").Value = "MATCH_CODE"

$battingLastRow = $batting.UsedRange.Rows.Count
$batting.Range("D2:D" + $battingLastRow).NumberFormat = "@"

for ($r = 2; $r -le $battingLastRow; $r++) {
    $link = $batting.Cells.Item($r, 4).Value2
    if ($link) {
        $code = $link -replace '.*MatchCode=', ''
        $batting.Cells.Item($r, 4).Value = $code
    }

    # Drop the INNING_NUMBER (column B) cell entirely when it carries no value
    $inning = $batting.Cells.Item($r, 2).Value2
    if (-not $inning) {
        $batting.Cells.Item($r, 2).ClearContents()
    }
}

# --------------------------------------------------------------------
# 3. "ODI Bowling" sheet updates
# --------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")

$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingLastRow = $bowling.UsedRange.Rows.Count
$bowling.Range("B2:B" + $bowlingLastRow).NumberFormat = "@"

for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $link = $bowling.Cells.Item($r, 2).Value2
    if ($link) {
        $code = $link -replace '.*MatchCode=', ''
        $bowling.Cells.Item($r, 2).Value = $code
    }
}

# --------------------------------------------------------------------
# 4. Append "ODI Batting Extra" sheet at the end
# --------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extraHeaders = New-Object 'object[,]' 1,6
$extraHeaders[0,0] = "MATCH_CODE"
$extraHeaders[0,1] = "BATTING_POSITION"
$extraHeaders[0,2] = "NUM_4"
$extraHeaders[0,3] = "NUM_6"
$extraHeaders[0,4] = "PERCENT_RUNS_OF_TOTAL"
$extraHeaders[0,5] = "MAN_OF_MATCH"
$extra.Range("A1:F1").Value = $extraHeaders
$extra.Range("A1:F1").Font.Bold = $true
$extra.Range("A1:F1").HorizontalAlignment = -4108
$extra.Range("A1:F1").VerticalAlignment = -4160
$extra.Range("A1:F1").Borders.LineStyle = 1

# MATCH_CODE (A) and MAN_OF_MATCH (F) are plain text, NUM_4/NUM_6/PERCENT
# columns (C/D/E) are also stored as text (including blanks); only
# BATTING_POSITION (B) is a real number when populated.
$extra.Range("A2:A21").NumberFormat = "@"
$extra.Range("C2:F21").NumberFormat = "@"

$extraRows = @(
    @("3942", 7,    "1", "0", "2.22%", "NO"),
    @("3983", 8,    "",  "",  "",      "NO"),
    @("3985", 8,    "0", "0", "0.65%", "NO"),
    @("3989", 7,    "0", "0", "0.27%", "NO"),
    @("3990", "",   "",  "",  "",      "NO"),
    @("3995", "",   "",  "",  "",      "NO"),
    @("3997", "",   "",  "",  "",      "NO"),
    @("3999", "",   "",  "",  "",      "NO"),
    @("4028", 8,    "2", "0", "7.12%", "NO"),
    @("4030", 8,    "",  "",  "",      "NO"),
    @("4033", 8,    "0", "0", "2.34%", "NO"),
    @("4037", "",   "",  "",  "",      "NO"),
    @("4517", "",   "",  "",  "",      "NO"),
    @("4557", 7,    "",  "",  "",      "NO"),
    @("4656", "",   "",  "",  "",      "NO"),
    @("4657", 7,    "0", "0", "5.76%", "NO"),
    @("4698", 7,    "0", "0", "0.67%", "NO"),
    @("4699", 8,    "",  "",  "",      "NO"),
    @("4700", 8,    "3", "1", "11.85%","NO"),
    @("4731", 8,    "1", "0", "1.52%", "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    $extra.Cells.Item($r, 3).Value = $row[2]
    $extra.Cells.Item($r, 4).Value = $row[3]
    $extra.Cells.Item($r, 5).Value = $row[4]
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}

$extra.Range("A1").Select()

# --------------------------------------------------------------------
# Activate the first sheet so the workbook opens on "Player Info"
# --------------------------------------------------------------------
$playerInfo.Activate()
